$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.547.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.560"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.52"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.300"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.079.62"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.32"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.825.37"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.561.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.13"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.64%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "172.56"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.85%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.41"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.71"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.60"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.419.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.674"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.80"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.955"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.75"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.59%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.978.73"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.58"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.10%  "

# Row 32 <-> Row 33 swap (Filecoin now in row 32, PancakeSwap now in row 33)
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.84"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.62%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.22%  "
